$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 02:16"

# Swap the Huesca / Huelva rows (province name + "Casos activos" value)
$ws.Range("A53").Value = "Huelva"
$ws.Range("C53").Value = 72

$ws.Range("A54").Value = "Huesca"
$ws.Range("C54").Value = 0
